$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.975.37"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.420.92"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'570.71"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "'140.03"
$ws.Range("E6").Value = "  -2.22%  "
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "2.404.77"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  -2.17%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "'26.22"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").Value = "2.851.60"
$ws.Range("D17").Value = "60.844.90"
$ws.Range("E17").Value = "  -2.12%  "
$ws.Range("D18").Value = "2.398.72"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = "  +6.79%  "
$ws.Range("D20").Value = "'10.66"
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "'323.19"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "'6.05"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").Value = "'64.83"
$ws.Range("E26").Value = "  -1.44%  "
$ws.Range("D27").Value = "'585.18"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("D28").Value = "'8.39"
$ws.Range("E28").Value = "  -9.07%  "
$ws.Range("D29").Value = "2.544.56"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("E32").Value = "  -4.57%  "
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  -5.60%  "
$ws.Range("D38").Value = "'151.82"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D40").Value = "'18.30"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "'5.14"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "'41.15"
$ws.Range("E44").Value = "  -3.97%  "
$ws.Range("D45").Value = "'2.37"
$ws.Range("E45").Value = "  -6.02%  "
$ws.Range("E46").Value = "  +13.04%  "
$ws.Range("D47").Value = "'142.22"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -3.48%  "
$ws.Range("D49").Value = "'0.590"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "'19.57"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("E51").Value = "  -3.53%  "
